$wb = $excel.ActiveWorkbook

# Duplicate the "FP" (1st/Final Pref) sheet as a template for the new
# "SP" (2nd Preference) sheet, placing it right after "FP".
$fp = $wb.Worksheets.Item("FP")
$fp.Copy([System.Reflection.Missing]::Value, $fp)
$sp = $wb.Worksheets.Item("FP (2)")
$sp.Name = "SP"

# New 2nd-preference source data (columns J:O, rows 2:5).
# Row 2 has no value for column M (left blank, as in the source data).
$sp.Range("J2").Value = 10.5
$sp.Range("K2").Value = 42.1
$sp.Range("L2").Value = 10.5
$sp.Range("M2").ClearContents()
$sp.Range("N2").Value = 26.3
$sp.Range("O2").Value = 10.5

$sp.Range("J3").Value = 5.3
$sp.Range("K3").Value = 26.3
$sp.Range("L3").Value = 10.5
$sp.Range("M3").Value = 0
$sp.Range("N3").Value = 50
$sp.Range("O3").Value = 7.9

$sp.Range("J4").Value = 3.7
$sp.Range("K4").Value = 18.5
$sp.Range("L4").Value = 11.1
$sp.Range("M4").Value = 0
$sp.Range("N4").Value = 55.6
$sp.Range("O4").Value = 11.1

$sp.Range("J5").Value = 3.2
$sp.Range("K5").Value = 15.9
$sp.Range("L5").Value = 11.1
$sp.Range("M5").Value = 0
$sp.Range("N5").Value = 60.3
$sp.Range("O5").Value = 9.5

# Selection/active-cell bookkeeping: FP is no longer the active tab, and
# its former single-cell selection is replaced by the full used range;
# the new SP sheet becomes the active tab with M2 selected.
[void]$fp.Range("A1:P5").Select()
[void]$sp.Range("M2").Select()
